# Update market-board derived values in each profession sheet
# (mirrors the scheduled runner's data refresh captured in the commit)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 287.7143
$ws.Range("I11").Value = 287.7143
$ws.Range("K11").Value = 287.7143
$ws.Range("M11").Value = -147.7143
$ws.Range("H12").Value = 448.8889
$ws.Range("I12").Value = 508
$ws.Range("K12").Value = 508
$ws.Range("M12").Value = -338
$ws.Range("H19").Value = 1757.625
$ws.Range("I19").Value = 2500
$ws.Range("K19").Value = 2500
$ws.Range("M19").Value = -2325
$ws.Range("H32").Value = 7098.6924
$ws.Range("J32").Value = 6598.778
$ws.Range("L32").Value = 6598.778
$ws.Range("N32").Value = -7250.778
$ws.Range("H40").Value = 8570.777
$ws.Range("I40").Value = 9694.5
$ws.Range("J40").Value = 8249.714
$ws.Range("K40").Value = 9694.5
$ws.Range("L40").Value = 8249.714
$ws.Range("M40").Value = -9519.5
$ws.Range("N40").Value = -8599.714
$ws.Range("H43").Value = 6902.353
$ws.Range("I43").Value = 7310
$ws.Range("J43").Value = 3845
$ws.Range("K43").Value = 7310
$ws.Range("L43").Value = 3845
$ws.Range("M43").Value = -7241
$ws.Range("N43").Value = -3983
$ws.Range("H51").Value = 11903.412
$ws.Range("I51").Value = 12443.556
$ws.Range("J51").Value = 11295.75
$ws.Range("K51").Value = 12443.556
$ws.Range("L51").Value = 11295.75
$ws.Range("M51").Value = -11959.556
$ws.Range("N51").Value = -12263.75
$ws.Range("H55").Value = 98.05882
$ws.Range("I55").Value = 43.42857
$ws.Range("J55").Value = 353
$ws.Range("K55").Value = 43.42857
$ws.Range("L55").Value = 353
$ws.Range("M55").Value = 170.57143
$ws.Range("N55").Value = -781
$ws.Range("H70").Value = 5485.8335
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 5485.8335
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 16457.5005
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -16997.5005
$ws.Range("H73").Value = 5485.8335
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 5485.8335
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 16457.5005
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -18329.5005
$ws.Range("H127").Value = 1172.8572
$ws.Range("I127").Value = 535
$ws.Range("K127").Value = 1605
$ws.Range("M127").Value = 3355
$ws.Range("H138").Value = 2393.6274
$ws.Range("I138").Value = 2160.5789
$ws.Range("J138").Value = 2532
$ws.Range("K138").Value = 6481.736699999999
$ws.Range("L138").Value = 7596
$ws.Range("M138").Value = -1341.736699999999
$ws.Range("N138").Value = -17876

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2361.652
$ws.Range("I45").Value = 1813.0588
$ws.Range("J45").Value = 3916
$ws.Range("K45").Value = 1813.0588
$ws.Range("L45").Value = 3916
$ws.Range("M45").Value = -1436.0588
$ws.Range("N45").Value = -4670
$ws.Range("H61").Value = 3898.25
$ws.Range("I61").Value = 2742.8076
$ws.Range("K61").Value = 2742.8076
$ws.Range("M61").Value = -2530.8076
$ws.Range("H132").Value = 1897.3214
$ws.Range("I132").Value = 1158.6538
$ws.Range("J132").Value = 11500
$ws.Range("K132").Value = 3475.9614
$ws.Range("L132").Value = 34500
$ws.Range("M132").Value = -945.9614000000001
$ws.Range("N132").Value = -39560
$ws.Range("H136").Value = 3898.25
$ws.Range("I136").Value = 2742.8076
$ws.Range("K136").Value = 8228.4228
$ws.Range("M136").Value = -5678.4228

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7388.95
$ws.Range("I20").Value = 1652.091
$ws.Range("J20").Value = 14400.667
$ws.Range("K20").Value = 1652.091
$ws.Range("L20").Value = 14400.667
$ws.Range("M20").Value = -1405.091
$ws.Range("N20").Value = -14894.667
$ws.Range("H80").Value = 57180.715
$ws.Range("J80").Value = 20052.8
$ws.Range("L80").Value = 20052.8
$ws.Range("N80").Value = -22048.8
$ws.Range("H83").Value = 57180.715
$ws.Range("J83").Value = 20052.8
$ws.Range("L83").Value = 100264
$ws.Range("N83").Value = -110248
$ws.Range("H134").Value = 4874.8857
$ws.Range("I134").Value = 3154.0667
$ws.Range("K134").Value = 9462.2001
$ws.Range("M134").Value = -6927.2001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3504.5293
$ws.Range("I22").Value = 1511.5454
$ws.Range("K22").Value = 1511.5454
$ws.Range("M22").Value = -1161.5454
$ws.Range("H38").Value = 39999
$ws.Range("J38").Value = 39999
$ws.Range("L38").Value = 39999
$ws.Range("N38").Value = -40753
$ws.Range("H39").Value = 21833
$ws.Range("I39").Value = 12750
$ws.Range("J39").Value = 39999
$ws.Range("K39").Value = 12750
$ws.Range("L39").Value = 39999
$ws.Range("M39").Value = -12359
$ws.Range("N39").Value = -40781
$ws.Range("H42").Value = 11659.333
$ws.Range("J42").Value = 23000
$ws.Range("L42").Value = 23000
$ws.Range("N42").Value = -24186
$ws.Range("H45").Value = 17249.75
$ws.Range("I45").Value = 11999.5
$ws.Range("K45").Value = 11999.5
$ws.Range("M45").Value = -11406.5
$ws.Range("H46").Value = 39999
$ws.Range("J46").Value = 39999
$ws.Range("L46").Value = 39999
$ws.Range("N46").Value = -40421
$ws.Range("H47").Value = 29999.5
$ws.Range("J47").Value = 29999.5
$ws.Range("L47").Value = 29999.5
$ws.Range("N47").Value = -31131.5
$ws.Range("H49").Value = 21833
$ws.Range("I49").Value = 12750
$ws.Range("J49").Value = 39999
$ws.Range("K49").Value = 12750
$ws.Range("L49").Value = 39999
$ws.Range("M49").Value = -12568
$ws.Range("N49").Value = -40363
$ws.Range("H54").Value = 33773.6
$ws.Range("J54").Value = 33773.6
$ws.Range("L54").Value = 33773.6
$ws.Range("N54").Value = -35089.6
$ws.Range("H56").Value = 14997
$ws.Range("J56").Value = 19994
$ws.Range("L56").Value = 19994
$ws.Range("N56").Value = -21684
$ws.Range("H59").Value = 105746.25
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1347
$ws.Range("I5").Value = 1194.3334
$ws.Range("K5").Value = 3583.0002
$ws.Range("M5").Value = -3471.0002
$ws.Range("H17").Value = 803.6667
$ws.Range("J17").Value = 650
$ws.Range("L17").Value = 1950
$ws.Range("N17").Value = -2288
$ws.Range("H39").Value = 8288.909
$ws.Range("I39").Value = 2750
$ws.Range("J39").Value = 9519.777
$ws.Range("K39").Value = 8250
$ws.Range("L39").Value = 28559.331
$ws.Range("M39").Value = -7956
$ws.Range("N39").Value = -29147.331
$ws.Range("H55").Value = 12249.25
$ws.Range("J55").Value = 12249.25
$ws.Range("L55").Value = 36747.75
$ws.Range("N55").Value = -37101.75
$ws.Range("H58").Value = 11929.667
$ws.Range("I58").Value = 11929.667
$ws.Range("K58").Value = 35789.001
$ws.Range("M58").Value = -35661.001
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("H134").Value = 2366.6667
$ws.Range("I134").Value = 2366.6667
$ws.Range("K134").Value = 7100.000100000001
$ws.Range("M134").Value = -2030.000100000001
$ws.Range("H135").Value = 1347
$ws.Range("I135").Value = 1194.3334
$ws.Range("K135").Value = 10749.0006
$ws.Range("M135").Value = -8214.000599999999
$ws.Range("H140").Value = 2075.3958
$ws.Range("I140").Value = 1576.1818
$ws.Range("K140").Value = 4728.5454
$ws.Range("M140").Value = 451.4546
$ws.Range("H141").Value = 7670.8823
$ws.Range("I141").Value = 6007
$ws.Range("K141").Value = 18021
$ws.Range("M141").Value = -12841

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5043.65
$ws.Range("I122").Value = 1742.0625
$ws.Range("J122").Value = 18250
$ws.Range("K122").Value = 5226.1875
$ws.Range("L122").Value = 54750
$ws.Range("M122").Value = -2776.1875
$ws.Range("N122").Value = -59650

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1597.2
$ws.Range("I22").Value = 1162.5
$ws.Range("K22").Value = 1162.5
$ws.Range("M22").Value = -867.5
$ws.Range("H27").Value = 1597.2
$ws.Range("I27").Value = 1162.5
$ws.Range("K27").Value = 1162.5
$ws.Range("M27").Value = -1055.5
$ws.Range("H46").Value = 1978.8334
$ws.Range("I46").Value = 1223.3334
$ws.Range("J46").Value = 2356.5833
$ws.Range("K46").Value = 1223.3334
$ws.Range("L46").Value = 2356.5833
$ws.Range("M46").Value = -1035.3334
$ws.Range("N46").Value = -2732.5833

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2402.9143
$ws.Range("I132").Value = 2056
$ws.Range("J132").Value = 2990
$ws.Range("K132").Value = 6168
$ws.Range("L132").Value = 8970
$ws.Range("M132").Value = -3638
$ws.Range("N132").Value = -14030
$ws.Range("H136").Value = 3825.9736
$ws.Range("I136").Value = 2611.1428
$ws.Range("J136").Value = 17999
$ws.Range("K136").Value = 7833.428400000001
$ws.Range("L136").Value = 53997
$ws.Range("M136").Value = -5283.428400000001
$ws.Range("N136").Value = -59097
